# Apply cryptos list update (commit: "Updated cryptos list on Fri Jun 14 14:36:44 UTC 2024 with GitHub Actions")
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'66.679.53"
$ws.Range("E2").Value = "'  -1.01%  "
$ws.Range("D3").Value = "'3.472.87"
$ws.Range("E3").Value = "'  -0.36%  "
$ws.Range("E4").Value = "'  -0.10%  "
$ws.Range("D5").Value = "'603.92"
$ws.Range("E5").Value = "'  -0.32%  "
$ws.Range("D6").Value = "'145.56"
$ws.Range("E6").Value = "'  -2.35%  "
$ws.Range("D7").Value = "'3.471.81"
$ws.Range("E7").Value = "'  -0.32%  "
$ws.Range("E8").Value = "'  -0.01%  "
$ws.Range("E9").Value = "'  -1.75%  "
$ws.Range("D10").Value = "'0.140"
$ws.Range("E10").Value = "'  -1.93%  "
$ws.Range("D11").Value = "'7.89"
$ws.Range("E11").Value = "'  +4.85%  "
$ws.Range("D12").Value = "'0.416"
$ws.Range("E12").Value = "'  -2.69%  "
$ws.Range("D13").Value = "'4.078.36"
$ws.Range("E13").Value = "'  +0.20%  "
$ws.Range("D14").Value = "'0.0000211"
$ws.Range("E14").Value = "'  -1.90%  "
$ws.Range("D15").Value = "'30.92"
$ws.Range("E15").Value = "'  -2.96%  "
$ws.Range("D16").Value = "'3.475.85"
$ws.Range("E16").Value = "'  -0.86%  "
$ws.Range("D17").Value = "'66.803.76"
$ws.Range("E17").Value = "'  -0.97%  "
$ws.Range("E18").Value = "'  -0.30%  "
$ws.Range("D19").Value = "'10.62"
$ws.Range("E19").Value = "'  +7.31%  "
$ws.Range("D20").Value = "'6.25"
$ws.Range("E20").Value = "'  -3.40%  "
$ws.Range("D21").Value = "'15.19"
$ws.Range("E21").Value = "'  -1.82%  "
$ws.Range("D22").Value = "'429.40"
$ws.Range("E22").Value = "'  -4.14%  "
$ws.Range("D23").Value = "'0.599"
$ws.Range("E23").Value = "'  -3.70%  "
$ws.Range("D24").Value = "'79.48"
$ws.Range("E24").Value = "'  +1.27%  "
$ws.Range("E25").Value = "'  +0.04%  "
$ws.Range("D26").Value = "'3.615.87"
$ws.Range("E26").Value = "'  -0.26%  "
$ws.Range("D27").Value = "'0.0000115"
$ws.Range("E27").Value = "'  -5.17%  "
$ws.Range("D28").Value = "'9.68"
$ws.Range("E28").Value = "'  -1.94%  "
$ws.Range("D29").Value = "'8.04"
$ws.Range("E29").Value = "'  -6.45%  "
$ws.Range("E30").Value = "'  -0.66%  "
$ws.Range("B31").Value = "'Binance-PegBSC-USD"
$ws.Range("C31").Value = "'https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd"
$ws.Range("D31").Value = "'1.00"
$ws.Range("E31").Value = "'  +0.05%  "
$ws.Range("B32").Value = "'Fetch.AI"
$ws.Range("C32").Value = "'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet"
$ws.Range("D32").Value = "'1.54"
$ws.Range("E32").Value = "'  -6.33%  "
$ws.Range("D33").Value = "'0.164"
$ws.Range("E33").Value = "'  -2.39%  "
$ws.Range("D34").Value = "'25.21"
$ws.Range("E34").Value = "'  -1.15%  "
$ws.Range("D35").Value = "'1.77"
$ws.Range("E35").Value = "'  -3.32%  "
$ws.Range("B36").Value = "'USDe"
$ws.Range("C36").Value = "'https://coinranking.com/coin/exbfr2U-0+usde-usde"
$ws.Range("D36").Value = "'1.00"
$ws.Range("E36").Value = "'  +0.02%  "
$ws.Range("B37").Value = "'NEARProtocol"
$ws.Range("C37").Value = "'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D37").Value = "'5.71"
$ws.Range("E37").Value = "'  -7.64%  "
$ws.Range("D38").Value = "'7.88"
$ws.Range("E38").Value = "'  -1.05%  "
$ws.Range("E39").Value = "'  +0.23%  "
$ws.Range("D40").Value = "'174.48"
$ws.Range("E40").Value = "'  -1.37%  "
$ws.Range("D41").Value = "'0.0886"
$ws.Range("E41").Value = "'  -1.01%  "
$ws.Range("D42").Value = "'5.31"
$ws.Range("E42").Value = "'  -1.49%  "
$ws.Range("D43").Value = "'1.98"
$ws.Range("E43").Value = "'  -13.00%  "
$ws.Range("D44").Value = "'0.886"
$ws.Range("E44").Value = "'  -1.09%  "
$ws.Range("D45").Value = "'46.29"
$ws.Range("E45").Value = "'  -0.52%  "
$ws.Range("D46").Value = "'27.52"
$ws.Range("E46").Value = "'  -10.63%  "
$ws.Range("D47").Value = "'1.20"
$ws.Range("E47").Value = "'  -6.73%  "
$ws.Range("D48").Value = "'7.27"
$ws.Range("E48").Value = "'  -4.19%  "
$ws.Range("D49").Value = "'2.37"
$ws.Range("E49").Value = "'  -3.97%  "
$ws.Range("D50").Value = "'0.972"
$ws.Range("E50").Value = "'  -2.27%  "
$ws.Range("D51").Value = "'0.243"
$ws.Range("E51").Value = "'  -2.42%  "
